$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Fix misplaced header cells (C16 -> B16, C25 -> B25)
$ws.Range("B16").Value2 = $ws.Range("C16").Value2
$ws.Range("C16").ClearContents()
$ws.Range("B25").Value2 = $ws.Range("C25").Value2
$ws.Range("C25").ClearContents()

# Chart 1: Lieu du stage
$co1 = $ws.ChartObjects().Add(0, 0, 0, 0)
$co1.Name = "Chart 1"
$chart1 = $co1.Chart
$chart1.ChartType = [Microsoft.Office.Interop.Excel.XlChartType]::xlPie
$chart1.SetSourceData($ws.Range("D10:E14"))
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Lieu du stage"

# Chart 2: Contenu du stage
$co2 = $ws.ChartObjects().Add(0, 0, 0, 0)
$co2.Name = "Chart 2"
$chart2 = $co2.Chart
$chart2.ChartType = [Microsoft.Office.Interop.Excel.XlChartType]::xlPie
$chart2.SetSourceData($ws.Range("D16:E23"))
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Contenu du stage"

# Chart 3: Type du stage
$co3 = $ws.ChartObjects().Add(0, 0, 0, 0)
$co3.Name = "Chart 3"
$chart3 = $co3.Chart
$chart3.ChartType = [Microsoft.Office.Interop.Excel.XlChartType]::xlPie
$chart3.SetSourceData($ws.Range("D25:E28"))
$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Type du stage"
